$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2392
$ws.Range("E2").Value = -6
$ws.Range("F2").Value = -6
$ws.Range("G2").Value = -22
$ws.Range("H2").Value = -22
$ws.Range("I2").Value = -22
$ws.Range("K2").Value = 1612
$ws.Range("L2").Value = 1047
$ws.Range("M2").Value = 565
$ws.Range("N2").Value = 565
$ws.Range("P2").Value = 187
$ws.Range("Q2").Value = 149
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = -238
$ws.Range("T2").Value = 83
$ws.Range("U2").Value = 66
$ws.Range("V2").Value = 490
$ws.Range("W2").Value = -0.24
$ws.Range("X2").Value = -0.9399999999999999
$ws.Range("Y2").Value = -5.97
$ws.Range("Z2").Value = -1.95
$ws.Range("AA2").Value = 185.14
$ws.Range("AB2").Value = 206.54
$ws.Range("AC2").Value = -64
$ws.Range("AD2").Value = -17.59
$ws.Range("AE2").Value = 1312
$ws.Range("AF2").Value = 0.86
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 43277382
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 3200
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 39
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 29
$ws.Range("K3").Value = 1643
$ws.Range("L3").Value = 1047
$ws.Range("M3").Value = 596
$ws.Range("N3").Value = 596
$ws.Range("P3").Value = 187
$ws.Range("Q3").Value = 81
$ws.Range("R3").Value = -70
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 127
$ws.Range("U3").Value = -46
$ws.Range("V3").Value = 494
$ws.Range("W3").Value = 0.2
$ws.Range("X3").Value = 0.9
$ws.Range("Y3").Value = 4.95
$ws.Range("Z3").Value = 1.77
$ws.Range("AA3").Value = 175.6
$ws.Range("AB3").Value = 221.85
$ws.Range("AC3").Value = 66
$ws.Range("AD3").Value = 25.85
$ws.Range("AE3").Value = 1383
$ws.Range("AF3").Value = 1.24
$ws.Range("AG3").Value = 43
$ws.Range("AH3").Value = 2.52
$ws.Range("AI3").Value = 10.82
$ws.Range("AJ3").Value = 43277382
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3270
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 17
$ws.Range("I4").Value = 17
$ws.Range("K4").Value = 1651
$ws.Range("L4").Value = 1062
$ws.Range("M4").Value = 588
$ws.Range("N4").Value = 588
$ws.Range("P4").Value = 187
$ws.Range("Q4").Value = 178
$ws.Range("R4").Value = -172
$ws.Range("S4").Value = -21
$ws.Range("T4").Value = 164
$ws.Range("U4").Value = 14
$ws.Range("V4").Value = 476
$ws.Range("W4").Value = 0.16
$ws.Range("X4").Value = 0.52
$ws.Range("Y4").Value = 2.89
$ws.Range("Z4").Value = 1.04
$ws.Range("AA4").Value = 180.54
$ws.Range("AB4").Value = 225.34
$ws.Range("AC4").Value = 40
$ws.Range("AD4").Value = 42.39
$ws.Range("AE4").Value = 1365
$ws.Range("AF4").Value = 1.23
$ws.Range("AG4").Value = 43
$ws.Range("AH4").Value = 2.58
$ws.Range("AI4").Value = 18.18
$ws.Range("AJ4").Value = 43277382
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3245
$ws.Range("E5").Value = -89
$ws.Range("F5").Value = -89
$ws.Range("G5").Value = -82
$ws.Range("H5").Value = -79
$ws.Range("I5").Value = -79
$ws.Range("K5").Value = 1783
$ws.Range("L5").Value = 1101
$ws.Range("M5").Value = 682
$ws.Range("N5").Value = 682
$ws.Range("P5").Value = 187
$ws.Range("Q5").Value = 46
$ws.Range("R5").Value = -21
$ws.Range("S5").Value = -23
$ws.Range("T5").Value = 115
$ws.Range("U5").Value = -69
$ws.Range("V5").Value = 456
$ws.Range("W5").Value = -2.73
$ws.Range("X5").Value = -2.44
$ws.Range("Y5").Value = -12.48
$ws.Range("Z5").Value = -4.62
$ws.Range("AA5").Value = 161.59
$ws.Range("AB5").Value = 183.26
$ws.Range("AC5").Value = -183
$ws.Range("AD5").Value = -7.06
$ws.Range("AE5").Value = 1581
$ws.Range("AF5").Value = 0.82
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 43277382
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 3358
$ws.Range("E6").Value = -160
$ws.Range("F6").Value = -160
$ws.Range("G6").Value = -440
$ws.Range("H6").Value = -441
$ws.Range("I6").Value = -441
$ws.Range("K6").Value = 1784
$ws.Range("L6").Value = 1562
$ws.Range("M6").Value = 221
$ws.Range("N6").Value = 221
$ws.Range("P6").Value = 187
$ws.Range("Q6").Value = -126
$ws.Range("R6").Value = -115
$ws.Range("S6").Value = 236
$ws.Range("T6").Value = 120
$ws.Range("U6").Value = -246
$ws.Range("V6").Value = 692
$ws.Range("W6").Value = -4.78
$ws.Range("X6").Value = -13.13
$ws.Range("Y6").Value = -97.63
$ws.Range("Z6").Value = -24.72
$ws.Range("AA6").Value = 705.72
$ws.Range("AB6").Value = -63.37
$ws.Range("AC6").Value = -1018
$ws.Range("AD6").Value = -2.18
$ws.Range("AE6").Value = 514
$ws.Range("AF6").Value = 4.32
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 43277382

# Row 7: clear financial data columns D:AI
$ws.Range("D7:AI7").ClearContents()

# Row 8: clear financial data columns D:AI
$ws.Range("D8:AI8").ClearContents()

# Row 9: clear financial data columns D:AI
$ws.Range("D9:AI9").ClearContents()
